# Auto-generated edit script replicating the OOXML diff for Sheets/Cerberus_Profits.xlsx
# (workbook tabs ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) — market-data refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("H13").Value = 11479.8
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 11479.8
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 11479.8
$ws.Range("N13").Value = -11817.8
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("H17").Value = 1969.4348
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1969.4348
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 5908.3044
$ws.Range("N17").Value = -6244.3044
$ws.Range("H74").Value = 5993.346
$ws.Range("I74").Value = 5628.75
$ws.Range("J74").Value = 6305.857
$ws.Range("K74").Value = 5628.75
$ws.Range("L74").Value = 6305.857
$ws.Range("M74").Value = -4692.75
$ws.Range("H77").Value = 5993.346
$ws.Range("I77").Value = 5628.75
$ws.Range("J77").Value = 6305.857
$ws.Range("K77").Value = 28143.75
$ws.Range("L77").Value = 31529.285
$ws.Range("M77").Value = -23463.75
$ws.Range("H92").Value = 5199
$ws.Range("I92").Value = 4999.25
$ws.Range("J92").Value = 5998
$ws.Range("K92").Value = 4999.25
$ws.Range("L92").Value = 5998
$ws.Range("M92").Value = -3751.25
$ws.Range("H96").Value = 777.6667
$ws.Range("I96").Value = 554.875
$ws.Range("J96").Value = 1223.25
$ws.Range("K96").Value = 1664.625
$ws.Range("L96").Value = 3669.75
$ws.Range("M96").Value = -291.625
$ws.Range("N96").Value = -6415.75
$ws.Range("H100").Value = 11037.777
$ws.Range("I100").Value = 3809.25
$ws.Range("J100").Value = 16820.6
$ws.Range("K100").Value = 3809.25
$ws.Range("L100").Value = 16820.6
$ws.Range("M100").Value = -3268.25
$ws.Range("N100").Value = -17902.6
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("H112").Value = 5832.524
$ws.Range("I112").Value = 1942
$ws.Range("J112").Value = 7777.7856
$ws.Range("K112").Value = 5826
$ws.Range("L112").Value = 23333.3568
$ws.Range("M112").Value = -4718
$ws.Range("N112").Value = -25549.3568
$ws.Range("H113").Value = 7223.55
$ws.Range("I113").Value = 6567.7144
$ws.Range("J113").Value = 8753.833000000001
$ws.Range("K113").Value = 6567.7144
$ws.Range("L113").Value = 8753.833000000001
$ws.Range("M113").Value = -3313.7144
$ws.Range("H131").Value = 3476.15
$ws.Range("I131").Value = 1460.1765
$ws.Range("J131").Value = 14900
$ws.Range("K131").Value = 4380.529500000001
$ws.Range("L131").Value = 44700
$ws.Range("M131").Value = 659.4704999999994
$ws.Range("H138").Value = 2936.8472
$ws.Range("I138").Value = 4945.4443
$ws.Range("J138").Value = 2267.3147
$ws.Range("K138").Value = 14836.3329
$ws.Range("L138").Value = 6801.9441
$ws.Range("M138").Value = -9696.332900000001
$ws.Range("N138").Value = -17081.9441
$ws.Range("N7").ClearContents()
$ws.Range("N14").ClearContents()
$ws.Range("N105").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1506.75
$ws.Range("I2").Value = 597.25
$ws.Range("J2").Value = 2416.25
$ws.Range("K2").Value = 597.25
$ws.Range("L2").Value = 2416.25
$ws.Range("M2").Value = -484.25
$ws.Range("N2").Value = -2642.25
$ws.Range("H22").Value = 9996
$ws.Range("I22").Value = 9996
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 9996
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -9697
$ws.Range("H74").Value = 1730.2439
$ws.Range("I74").Value = 927.63635
$ws.Range("J74").Value = 2659.5789
$ws.Range("K74").Value = 927.63635
$ws.Range("L74").Value = 2659.5789
$ws.Range("M74").Value = -53.63634999999999
$ws.Range("N74").Value = -4407.5789
$ws.Range("H77").Value = 1730.2439
$ws.Range("I77").Value = 927.63635
$ws.Range("J77").Value = 2659.5789
$ws.Range("K77").Value = 4638.18175
$ws.Range("L77").Value = 13297.8945
$ws.Range("M77").Value = -270.1817499999997
$ws.Range("N77").Value = -22033.8945
$ws.Range("H102").Value = 3747.4375
$ws.Range("I102").Value = 3709.7856
$ws.Range("J102").Value = 4011
$ws.Range("K102").Value = 3709.7856
$ws.Range("L102").Value = 4011
$ws.Range("M102").Value = -2087.7856
$ws.Range("H116").Value = 1506.75
$ws.Range("I116").Value = 597.25
$ws.Range("J116").Value = 2416.25
$ws.Range("K116").Value = 597.25
$ws.Range("L116").Value = 2416.25
$ws.Range("M116").Value = 1696.75
$ws.Range("N116").Value = -7004.25
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1506.75
$ws.Range("I3").Value = 597.25
$ws.Range("J3").Value = 2416.25
$ws.Range("K3").Value = 597.25
$ws.Range("L3").Value = 2416.25
$ws.Range("M3").Value = -483.25
$ws.Range("N3").Value = -2644.25
$ws.Range("H20").Value = 2190.6667
$ws.Range("I20").Value = 2190.6667
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 2190.6667
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -1943.6667
$ws.Range("H94").Value = 9448.1
$ws.Range("I94").Value = 5032.577
$ws.Range("J94").Value = 17648.357
$ws.Range("K94").Value = 5032.577
$ws.Range("L94").Value = 17648.357
$ws.Range("M94").Value = -4581.577
$ws.Range("N20").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1623.4166
$ws.Range("I58").Value = 1314.8334
$ws.Range("J58").Value = 1932
$ws.Range("K58").Value = 1314.8334
$ws.Range("L58").Value = 1932
$ws.Range("M58").Value = -1111.8334
$ws.Range("N58").Value = -2338
$ws.Range("H86").Value = 5037.8335
$ws.Range("I86").Value = 4541.5
$ws.Range("J86").Value = 5534.1665
$ws.Range("K86").Value = 4541.5
$ws.Range("L86").Value = 5534.1665
$ws.Range("M86").Value = -3418.5
$ws.Range("H89").Value = 5037.8335
$ws.Range("I89").Value = 4541.5
$ws.Range("J89").Value = 5534.1665
$ws.Range("K89").Value = 22707.5
$ws.Range("L89").Value = 27670.8325
$ws.Range("M89").Value = -17091.5
$ws.Range("H134").Value = 1589.3864
$ws.Range("I134").Value = 1545.258
$ws.Range("J134").Value = 1694.6154
$ws.Range("K134").Value = 4635.774
$ws.Range("L134").Value = 5083.8462
$ws.Range("M134").Value = -2100.774
$ws.Range("H136").Value = 1623.4166
$ws.Range("I136").Value = 1314.8334
$ws.Range("J136").Value = 1932
$ws.Range("K136").Value = 3944.5002
$ws.Range("L136").Value = 5796
$ws.Range("M136").Value = -1394.5002
$ws.Range("N136").Value = -10896

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 10640520
$ws.Range("I129").Value = 23810304
$ws.Range("J129").Value = 8335808
$ws.Range("K129").Value = 71430912
$ws.Range("L129").Value = 25007424
$ws.Range("M129").Value = -71425912
$ws.Range("N129").Value = -25017424
$ws.Range("H131").Value = 27779746
$ws.Range("I131").Value = 22223912
$ws.Range("J131").Value = 30305126
$ws.Range("K131").Value = 66671736
$ws.Range("L131").Value = 90915378
$ws.Range("M131").Value = -66666696
$ws.Range("N131").Value = -90925458

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 4499.6665
$ws.Range("I22").Value = 1749.5
$ws.Range("J22").Value = 10000
$ws.Range("K22").Value = 1749.5
$ws.Range("L22").Value = 10000
$ws.Range("M22").Value = -1220.5
$ws.Range("N22").Value = -11058
$ws.Range("H102").Value = 15983.429
$ws.Range("I102").Value = 27758.6
$ws.Range("J102").Value = 5278.727
$ws.Range("K102").Value = 27758.6
$ws.Range("L102").Value = 5278.727
$ws.Range("M102").Value = -26136.6
$ws.Range("N102").Value = -8522.726999999999
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 557.6316
$ws.Range("I16").Value = 406.17648
$ws.Range("J16").Value = 1845
$ws.Range("K16").Value = 406.17648
$ws.Range("L16").Value = 1845
$ws.Range("M16").Value = -236.17648
$ws.Range("H40").Value = 2642.875
$ws.Range("I40").Value = 2199
$ws.Range("J40").Value = 3086.75
$ws.Range("K40").Value = 2199
$ws.Range("L40").Value = 3086.75
$ws.Range("M40").Value = -2063
$ws.Range("N40").Value = -3358.75
$ws.Range("H68").Value = 2458.7222
$ws.Range("I68").Value = 2437.2
$ws.Range("J68").Value = 2566.3333
$ws.Range("K68").Value = 2437.2
$ws.Range("L68").Value = 2566.3333
$ws.Range("M68").Value = -1688.2
$ws.Range("N68").Value = -4064.3333
$ws.Range("H71").Value = 2458.7222
$ws.Range("I71").Value = 2437.2
$ws.Range("J71").Value = 2566.3333
$ws.Range("K71").Value = 12186
$ws.Range("L71").Value = 12831.6665
$ws.Range("M71").Value = -8442
$ws.Range("N71").Value = -20319.6665
$ws.Range("H97").Value = 20936.6
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 20936.6
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 20936.6
$ws.Range("N97").Value = -22918.6
$ws.Range("H108").Value = 49500
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 49500
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 49500
$ws.Range("N108").Value = -57180
$ws.Range("H109").Value = 64284.5
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 64284.5
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 64284.5
$ws.Range("N109").Value = -67058.5
$ws.Range("H132").Value = 2625.1155
$ws.Range("I132").Value = 2539.5
$ws.Range("J132").Value = 2640.682
$ws.Range("K132").Value = 7618.5
$ws.Range("L132").Value = 7922.045999999999
$ws.Range("M132").Value = -5088.5
$ws.Range("N132").Value = -12982.046
$ws.Range("H136").Value = 2146.2666
$ws.Range("I136").Value = 2046.4231
$ws.Range("J136").Value = 2282.8948
$ws.Range("K136").Value = 6139.2693
$ws.Range("L136").Value = 6848.6844
$ws.Range("M136").Value = -3589.2693
$ws.Range("N136").Value = -11948.6844

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 5499
$ws.Range("I14").Value = 5748.5
$ws.Range("J14").Value = 5000
$ws.Range("K14").Value = 5748.5
$ws.Range("L14").Value = 5000
$ws.Range("M14").Value = -5580.5
$ws.Range("H81").Value = 9316.333000000001
$ws.Range("I81").Value = 9316.333000000001
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 18632.666
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -17571.666
$ws.Range("H84").Value = 9316.333000000001
$ws.Range("I84").Value = 9316.333000000001
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 93163.33
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -87859.33
$ws.Range("H109").Value = 69376.5
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 69376.5
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 69376.5
$ws.Range("N109").Value = -72150.5
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("H126").Value = 2321.1538
$ws.Range("I126").Value = 2281.3333
$ws.Range("J126").Value = 2799
$ws.Range("K126").Value = 6843.999899999999
$ws.Range("L126").Value = 8397
$ws.Range("M126").Value = -4373.999899999999
$ws.Range("N126").Value = -13337
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
